# -----------------------------------------------------------------------
# Applies the "Summary" section, certification heading font-size bump,
# mid-word bookmark relocation ("_GoBack"), and the three new
# Unity/Report/Engage sections described by the commit diff.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# 1) After the title "Individual Critical Appraisal", add a blank paragraph
#    followed by a new "Summary " Heading 2 paragraph.
$pTitle = $d.Paragraphs(1)
$rTitle = $pTitle.Range
$rTitle.Collapse(0)
$rTitle.InsertAfter("`r")

$pBlank = $d.Paragraphs(2)
$pBlank.Range.InsertParagraphAfter()
$pSummary = $d.Paragraphs(3)
$pSummary.Style = "Heading 2"
$pSummary.Range.Text = "Summary "

# 2) Bump the "My contributions" Heading 2 run/paragraph mark to 12pt (sz 24).
$pMyContrib = $d.Paragraphs(5)
$pMyContrib.Range.Font.Size = 12

# 3) Split the "Providing a large chunk..." run so a relocated "_GoBack"
#    bookmark sits after "...source cod".
$pContribution = $d.Paragraphs(6)
$splitOffset = $pContribution.Range.Start + 41
$splitRange = $d.Range($splitOffset, $splitOffset)
$d.Bookmarks.Add("_GoBack", $splitRange)

# 4) Bump the "My Failings" Heading 2 run/paragraph mark to 12pt (sz 24).
$pMyFailings = $d.Paragraphs(9)
$pMyFailings.Range.Font.Size = 12

# 5) Append the three new "Unity Work" / "Report Work" / "Engage Work"
#    sections after the final existing paragraph.
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$pLast.Range.InsertParagraphAfter()
$pUnityHeading = $d.Paragraphs($d.Paragraphs.Count)
$pUnityHeading.Style = "Heading 2"
$pUnityHeading.Range.Text = "Unity Work"

$pUnityHeading.Range.InsertParagraphAfter()
$pUnityBody = $d.Paragraphs($d.Paragraphs.Count)
$pUnityBody.Style = "Normal"
$pUnityBody.Range.Text = "List unity work in more detail"

$pUnityBody.Range.InsertParagraphAfter()
$pReportHeading = $d.Paragraphs($d.Paragraphs.Count)
$pReportHeading.Style = "Heading 2"
$pReportHeading.Range.Text = "Report Work"

$pReportHeading.Range.InsertParagraphAfter()
$pReportBody = $d.Paragraphs($d.Paragraphs.Count)
$pReportBody.Style = "Normal"
$pReportBody.Range.Text = "List report work in more detail"

$pReportBody.Range.InsertParagraphAfter()
$pEngageHeading = $d.Paragraphs($d.Paragraphs.Count)
$pEngageHeading.Style = "Heading 2"
$pEngageHeading.Range.Text = "Engage Work"

$pEngageHeading.Range.InsertParagraphAfter()
$pEngageBody = $d.Paragraphs($d.Paragraphs.Count)
$pEngageBody.Style = "Normal"
$pEngageBody.Range.Text = "List engage work in more detail"
